$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that follows the title heading.
$metaPara = $d.Paragraphs(2)
if ($metaPara.Range.Text -match "^Meta description") {
    $metaPara.Range.Delete()
}

# 2. Insert a new bold paragraph ("Play Easter Island free | Review of Easter
#    Island Online Slot") right before the final (italic) paragraph.
$count = $d.Paragraphs.Count
$precedingPara = $d.Paragraphs($count - 1)
$precedingPara.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$newPara.Range.Style = "Normal"

$fullRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Easter Island free | Review of Easter Island Online Slot</w:t></w:r></w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$fullRange.InsertXML($xml) | Out-Null

# 3. Replace the text of the final italic paragraph (feature-image prompt ->
#    meta description text), keeping its italic formatting intact.
$d.Content.Find.Execute(
    "Create a cartoon-style feature image for Easter Island that showcases a happy Maya warrior with glasses. The image should have a fun and vibrant feel, with the warrior holding a bunny or Easter egg to tie in with the Easter theme of the game. Use a colorful and engaging background, such as palm trees on a beach or a field of flowers, to make the image stand out. Make sure the Maya warrior has a big smile on their face to convey the fun and exciting atmosphere of the game.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Easter Island slot game. Play Easter Island free & enjoy the Prize Pick and Free Spins bonus features. Easter-themed symbols, interesting graphics and sound.",
    2) | Out-Null
